# Adicionar especificações atribuídas ao Cesário e ao Tiago
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The exception text for "Pintor já existe no sistema" (step 2) is refined
$ws.Range("A16").Value = "Excepção 1               (passo 2)`n[Pintor já existe no sistema]"

# Move the active selection to A20, matching the saved view state
$ws.Range("A20").Select() | Out-Null
